$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression (only B2 changes slightly)
$ws.Range("B2").Value = 0.2240730911947776

# Row 3 - RandomForestRegressor (values change, name stays the same)
$ws.Range("B3").Value = 0.1657081008439971
$ws.Range("C3").Value = 0.1552878192046485
$ws.Range("D3").Value = 0.1662867392806648

# Row 4 - model renamed from GradientBoostingRegressor to DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.1820579342582942
$ws.Range("C4").Value = 0.1826762870675177
$ws.Range("D4").Value = 0.1672690861103096

# Row 5 - model renamed from AdaBoostRegressor to MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.167039699821732
$ws.Range("C5").Value = 0.1408758275932796
$ws.Range("D5").Value = 0.1764768898534349
